$wb = $excel.ActiveWorkbook

$jan = $wb.Worksheets.Item("January")
$feb = $wb.Worksheets.Item("February")

# "Updated tasks for Feb": copy January's Name/Task table into February.
$jan.Range("B2:C5").Copy() | Out-Null
$feb.Range("B2").PasteSpecial() | Out-Null
$excel.CutCopyMode = $false

# Column C needs to be wide enough to fit the task names on both sheets.
$jan.Columns.Item(3).ColumnWidth = 18.6
$feb.Columns.Item(3).ColumnWidth = 18.6

# Selections: January cursor moves to C11, February cursor sits on C5.
$jan.Range("C11").Select() | Out-Null
$feb.Range("C5").Select() | Out-Null

# February becomes the active/visible tab.
$feb.Activate() | Out-Null
